# Update the answer key table: "two-digit number divided by one-digit
# number" practice sheet. Each of the 25 populated table cells gets its
# division problem replaced with a new one (text-only change; the table
# shape/cell count is unchanged).

$d = $word.ActiveDocument

$pairs = @(
    @{ Old = "60÷5=12, 0";  New = "98÷2=49, 0" },
    @{ Old = "81÷7=11, 4";  New = "42÷4=10, 2" },
    @{ Old = "92÷9=10, 2";  New = "14÷3=4, 2" },
    @{ Old = "52÷7=7, 3";   New = "20÷3=6, 2" },
    @{ Old = "68÷4=17, 0";  New = "55÷4=13, 3" },
    @{ Old = "72÷9=8, 0";   New = "49÷9=5, 4" },
    @{ Old = "65÷6=10, 5";  New = "17÷2=8, 1" },
    @{ Old = "71÷2=35, 1";  New = "62÷3=20, 2" },
    @{ Old = "15÷8=1, 7";   New = "68÷4=17, 0" },
    @{ Old = "13÷7=1, 6";   New = "43÷3=14, 1" },
    @{ Old = "30÷3=10, 0";  New = "17÷5=3, 2" },
    @{ Old = "58÷8=7, 2";   New = "15÷7=2, 1" },
    @{ Old = "62÷2=31, 0";  New = "63÷5=12, 3" },
    @{ Old = "83÷4=20, 3";  New = "23÷5=4, 3" },
    @{ Old = "33÷9=3, 6";   New = "32÷7=4, 4" },
    @{ Old = "86÷5=17, 1";  New = "44÷5=8, 4" },
    @{ Old = "17÷7=2, 3";   New = "89÷5=17, 4" },
    @{ Old = "59÷4=14, 3";  New = "73÷3=24, 1" },
    @{ Old = "57÷7=8, 1";   New = "90÷5=18, 0" },
    @{ Old = "12÷4=3, 0";   New = "74÷4=18, 2" },
    @{ Old = "61÷8=7, 5";   New = "29÷5=5, 4" },
    @{ Old = "98÷9=10, 8";  New = "81÷6=13, 3" },
    @{ Old = "63÷7=9, 0";   New = "55÷2=27, 1" },
    @{ Old = "10÷4=2, 2";   New = "10÷8=1, 2" },
    @{ Old = "96÷8=12, 0";  New = "98÷9=10, 8" }
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.New, 2)
}
